$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update policy configuration values (row 2)
$ws.Range("B2").Value = "scope:Absent"
$ws.Range("D2").Value = "automation"

# Move the active cell selection from D2 to D3
$ws.Range("D3").Select()
